$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1581.909
$ws.Range("I40").Value = 1750.125
$ws.Range("K40").Value = 1750.125
$ws.Range("M40").Value = -1575.125
$ws.Range("H64").Value = 6173.5
$ws.Range("I64").Value = 6338.8
$ws.Range("K64").Value = 6338.8
$ws.Range("M64").Value = -6090.8
$ws.Range("H67").Value = 6173.5
$ws.Range("I67").Value = 6338.8
$ws.Range("K67").Value = 6338.8
$ws.Range("M67").Value = -5480.8
$ws.Range("H98").Value = 1460.862
$ws.Range("I98").Value = 1460.862
$ws.Range("K98").Value = 1460.862
$ws.Range("M98").Value = 37.13799999999992
$ws.Range("H113").Value = 2977.6667
$ws.Range("J113").Value = 3185.5715
$ws.Range("L113").Value = 3185.5715
$ws.Range("N113").Value = -9693.5715
$ws.Range("H122").Value = 1460.862
$ws.Range("I122").Value = 1460.862
$ws.Range("K122").Value = 4382.586
$ws.Range("M122").Value = -1932.586
$ws.Range("H138").Value = 2350.0588
$ws.Range("I138").Value = 1703.6
$ws.Range("J138").Value = 3273.5715
$ws.Range("K138").Value = 5110.799999999999
$ws.Range("L138").Value = 9820.7145
$ws.Range("M138").Value = 29.20000000000073
$ws.Range("N138").Value = -20100.7145

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 43707.117
$ws.Range("I32").Value = 49103.695
$ws.Range("K32").Value = 49103.695
$ws.Range("M32").Value = -48816.695
$ws.Range("H61").Value = 4375.6924
$ws.Range("I61").Value = 3303.0557
$ws.Range("K61").Value = 3303.0557
$ws.Range("M61").Value = -3091.0557
$ws.Range("H74").Value = 2526.3784
$ws.Range("I74").Value = 523.1177
$ws.Range("J74").Value = 4229.15
$ws.Range("K74").Value = 523.1177
$ws.Range("L74").Value = 4229.15
$ws.Range("M74").Value = 350.8823
$ws.Range("N74").Value = -5977.15
$ws.Range("H77").Value = 2526.3784
$ws.Range("I77").Value = 523.1177
$ws.Range("J77").Value = 4229.15
$ws.Range("K77").Value = 2615.5885
$ws.Range("L77").Value = 21145.75
$ws.Range("M77").Value = 1752.4115
$ws.Range("N77").Value = -29881.75
$ws.Range("H110").Value = 29690198
$ws.Range("I110").Value = 45240068
$ws.Range("K110").Value = 45240068
$ws.Range("M110").Value = -45238023
$ws.Range("H136").Value = 4375.6924
$ws.Range("I136").Value = 3303.0557
$ws.Range("K136").Value = 9909.167099999999
$ws.Range("M136").Value = -7359.167099999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 4085
$ws.Range("I20").Value = 3404.2856
$ws.Range("J20").Value = 4879.1665
$ws.Range("K20").Value = 3404.2856
$ws.Range("L20").Value = 4879.1665
$ws.Range("M20").Value = -3157.2856
$ws.Range("N20").Value = -5373.1665
$ws.Range("H86").Value = 290390.84
$ws.Range("I86").Value = 5456
$ws.Range("K86").Value = 5456
$ws.Range("M86").Value = -4333
$ws.Range("H89").Value = 290390.84
$ws.Range("I89").Value = 5456
$ws.Range("K89").Value = 27280
$ws.Range("M89").Value = -21664
$ws.Range("H94").Value = 1265.1333
$ws.Range("I94").Value = 958.16
$ws.Range("K94").Value = 958.16
$ws.Range("M94").Value = -507.16
$ws.Range("H134").Value = 3531.3928
$ws.Range("I134").Value = 1560.2858
$ws.Range("K134").Value = 4680.857400000001
$ws.Range("M134").Value = -2145.857400000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 5004.1333
$ws.Range("I16").Value = 4155
$ws.Range("K16").Value = 4155
$ws.Range("M16").Value = -3868
$ws.Range("H94").Value = 16333.333
$ws.Range("I94").Value = 23750
$ws.Range("J94").Value = 1500
$ws.Range("K94").Value = 23750
$ws.Range("L94").Value = 1500
$ws.Range("M94").Value = -23299
$ws.Range("N94").Value = -2402
$ws.Range("H113").Value = 5004.1333
$ws.Range("I113").Value = 4155
$ws.Range("K113").Value = 4155
$ws.Range("M113").Value = -1985
$ws.Range("H122").Value = 126167.5
$ws.Range("I122").Value = 167890.17
$ws.Range("K122").Value = 503670.51
$ws.Range("M122").Value = -501220.51

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H18").Value = 555.5
$ws.Range("I18").Value = 572.7778
$ws.Range("K18").Value = 1718.3334
$ws.Range("M18").Value = -1549.3334
$ws.Range("H47").Value = 111819.336
$ws.Range("I47").Value = 143266.28
$ws.Range("K47").Value = 429798.84
$ws.Range("M47").Value = -429367.84
$ws.Range("H51").Value = 3389.25
$ws.Range("I51").Value = 302
$ws.Range("J51").Value = 4418.3335
$ws.Range("K51").Value = 906
$ws.Range("L51").Value = 13255.0005
$ws.Range("M51").Value = -446
$ws.Range("N51").Value = -14175.0005
$ws.Range("H55").Value = 5187
$ws.Range("I55").Value = 375
$ws.Range("J55").Value = 9999
$ws.Range("K55").Value = 1125
$ws.Range("L55").Value = 29997
$ws.Range("M55").Value = -948
$ws.Range("N55").Value = -30351
$ws.Range("H111").Value = 14375
$ws.Range("I111").Value = 500
$ws.Range("K111").Value = 1500
$ws.Range("M111").Value = 1567
$ws.Range("H112").Value = 2510579.8
$ws.Range("I112").Value = 6667666
$ws.Range("J112").Value = 16328
$ws.Range("K112").Value = 20002998
$ws.Range("L112").Value = 48984
$ws.Range("M112").Value = -20001890
$ws.Range("N112").Value = -51200

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 20830.5
$ws.Range("I70").Value = 18861
$ws.Range("J70").Value = 24113
$ws.Range("K70").Value = 18861
$ws.Range("L70").Value = 24113
$ws.Range("M70").Value = -18591
$ws.Range("N70").Value = -24653
$ws.Range("H73").Value = 20830.5
$ws.Range("I73").Value = 18861
$ws.Range("J73").Value = 24113
$ws.Range("K73").Value = 18861
$ws.Range("L73").Value = 24113
$ws.Range("M73").Value = -17925
$ws.Range("N73").Value = -25985
$ws.Range("H100").Value = 14224.75
$ws.Range("J100").Value = 14224.75
$ws.Range("L100").Value = 14224.75
$ws.Range("N100").Value = -16388.75
$ws.Range("H102").Value = 2649.2778
$ws.Range("I102").Value = 2477.8462
$ws.Range("K102").Value = 2477.8462
$ws.Range("M102").Value = -855.8462
$ws.Range("H121").Value = 0
$ws.Range("J121").Value = 0
$ws.Range("L121").ClearContents()
$ws.Range("N121").Value = 0
$ws.Range("H126").Value = 2536.2354
$ws.Range("I126").Value = 2374.923
$ws.Range("K126").Value = 7124.768999999999
$ws.Range("M126").Value = -4654.768999999999
$ws.Range("H132").Value = 4299.353
$ws.Range("I132").Value = 2196.8948
$ws.Range("K132").Value = 6590.6844
$ws.Range("M132").Value = -4060.6844
$ws.Range("H136").Value = 121081.25
$ws.Range("J136").Value = 121081.25
$ws.Range("L136").Value = 363243.75
$ws.Range("N136").Value = -368343.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H3").Value = 5000
$ws.Range("J3").Value = 5000
$ws.Range("L3").Value = 5000
$ws.Range("N3").Value = -5224
$ws.Range("H7").Value = 9521.182000000001
$ws.Range("I7").Value = 6588.8335
$ws.Range("J7").Value = 13040
$ws.Range("K7").Value = 6588.8335
$ws.Range("L7").Value = 13040
$ws.Range("M7").Value = -6476.8335
$ws.Range("N7").Value = -13264
$ws.Range("H15").Value = 5000
$ws.Range("J15").Value = 5000
$ws.Range("L15").Value = 5000
$ws.Range("N15").Value = -5340
$ws.Range("H22").Value = 2930.4285
$ws.Range("I22").Value = 1792.52
$ws.Range("J22").Value = 4603.8237
$ws.Range("K22").Value = 1792.52
$ws.Range("L22").Value = 4603.8237
$ws.Range("M22").Value = -1497.52
$ws.Range("N22").Value = -5193.8237
$ws.Range("H24").Value = 10000
$ws.Range("J24").Value = 10000
$ws.Range("L24").Value = 10000
$ws.Range("N24").Value = -10686
$ws.Range("H27").Value = 2930.4285
$ws.Range("I27").Value = 1792.52
$ws.Range("J27").Value = 4603.8237
$ws.Range("K27").Value = 1792.52
$ws.Range("L27").Value = 4603.8237
$ws.Range("M27").Value = -1685.52
$ws.Range("N27").Value = -4817.8237
$ws.Range("H40").Value = 14451.706
$ws.Range("I40").Value = 13695.357
$ws.Range("J40").Value = 17981.334
$ws.Range("K40").Value = 13695.357
$ws.Range("L40").Value = 17981.334
$ws.Range("M40").Value = -13559.357
$ws.Range("N40").Value = -18253.334
$ws.Range("H55").Value = 349.375
$ws.Range("I55").Value = 300.53845
$ws.Range("K55").Value = 300.53845
$ws.Range("M55").Value = -127.53845
$ws.Range("H126").Value = 9521.182000000001
$ws.Range("I126").Value = 6588.8335
$ws.Range("J126").Value = 13040
$ws.Range("K126").Value = 19766.5005
$ws.Range("L126").Value = 39120
$ws.Range("M126").Value = -17296.5005
$ws.Range("N126").Value = -44060

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H30").Value = 50000
$ws.Range("J30").Value = 50000
$ws.Range("L30").Value = 50000
$ws.Range("N30").Value = -50214
$ws.Range("H96").Value = 41679
$ws.Range("I96").Value = 2796.6667
$ws.Range("J96").Value = 100002.5
$ws.Range("K96").Value = 2796.6667
$ws.Range("L96").Value = 100002.5
$ws.Range("M96").Value = -1423.6667
$ws.Range("N96").Value = -102748.5
$ws.Range("H132").Value = 5362.636
$ws.Range("I132").Value = 2761.25
$ws.Range("J132").Value = 7811
$ws.Range("K132").Value = 8283.75
$ws.Range("L132").Value = 23433
$ws.Range("M132").Value = -5753.75
$ws.Range("N132").Value = -28493
$ws.Range("H138").Value = 85000
$ws.Range("J138").Value = 85000
$ws.Range("L138").Value = 85000
$ws.Range("N138").Value = -95280
